$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.855.04"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.55"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.81"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.871.25"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.641.35"
$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.48"
$ws.Range("E16").Value = "  +1.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.865.93"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.41"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  +6.63%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.66"
$ws.Range("E25").Value = "  +1.85%  "

$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("E28").Value = "  +1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.76"
$ws.Range("E29").Value = "  +0.67%  "

$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("E32").Value = "  +1.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.281.73"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +0.60%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.819"
$ws.Range("E39").Value = "  -0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.783.12"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  -6.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.68"
$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.03"
$ws.Range("E46").Value = "  -0.57%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0967"
$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("E51").Value = "  -0.21%  "
